# Adds the 5 new "Royal Glory" / "Springtime" Durazno price rows recorded
# for the Macroferia Regional de Talca on 2022-12-16 (serial 44911),
# appending them right after the last existing data row (row 438).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 439

$newRows = @(
    @(5, "Macroferia Regional de Talca", "Maule", 44911, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Royal Glory",  "Primera", 250, 12000, 12000, 12000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 800, 15),
    @(5, "Macroferia Regional de Talca", "Maule", 44911, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Royal Glory",  "Segunda", 200, 10000, 10000, 10000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 667, 15),
    @(5, "Macroferia Regional de Talca", "Maule", 44911, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Springtime",  "Especial", 180, 14000, 14000, 14000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 933, 15),
    @(5, "Macroferia Regional de Talca", "Maule", 44911, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Springtime",  "Primera", 150, 12000, 12000, 12000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 800, 15),
    @(5, "Macroferia Regional de Talca", "Maule", 44911, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Springtime",  "Segunda", 120, 10000, 10000, 10000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 667, 15)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
    # Column D (4) holds a date serial; give it the same date number format
    # used by every other row in the sheet (style index 2 / "YYYY-MM-DD HH:MM:SS").
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($startRow - 1, 4).NumberFormat()
}
